# Doing Updates for Financials
# Insert two new columns (D:E) before the existing quarterly data so a new
# quarter (Mar-19) and the prior quarter (Dec-18) can be added in front of
# the existing series, which shifts right by two columns (old D:K -> F:M).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HOS")

# --- 1. Insert two blank columns at D:E, pushing old D:K to F:M -----------
$ws.Range("D1:E1").EntireColumn.Insert()

# --- 2. Restore number formatting on the freshly inserted D:E columns -----
# EntireColumn.Insert() pulls formatting from the column to the left (C),
# but this sheet uses a uniform per-row style across D:K (style 2 for the
# "Period Ending" date rows, style 3 for the numeric data rows). Paste the
# format from column F (which kept its original per-row style) back onto
# the new D:E cells so no new style entries are introduced.
$ws.Range("F7").Copy()
$ws.Range("D7:E7").PasteSpecial(-4122)
$ws.Range("F38").Copy()
$ws.Range("D38:E38").PasteSpecial(-4122)
$ws.Range("F80").Copy()
$ws.Range("D80:E80").PasteSpecial(-4122)

$ws.Range("F8").Copy()
$ws.Range("D8:E35").PasteSpecial(-4122)
$ws.Range("F39").Copy()
$ws.Range("D39:E77").PasteSpecial(-4122)
$ws.Range("F81").Copy()
$ws.Range("D81:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 3. Populate the new D:E columns with the new quarter's figures -------
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 53900
$ws.Range("E8").Value = 58500
$ws.Range("D9").Value = 38600
$ws.Range("E9").Value = 38200
$ws.Range("D10").Value = 15300
$ws.Range("E10").Value = 20300
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 27600
$ws.Range("E15").Value = 27600
$ws.Range("D17").Value = 69500
$ws.Range("E17").Value = 80900
$ws.Range("D18").Value = -15600
$ws.Range("E18").Value = -22400
$ws.Range("D20").Value = 600
$ws.Range("E20").Value = 500
$ws.Range("D21").Value = 12600
$ws.Range("E21").Value = 5700
$ws.Range("D22").Value = 16700
$ws.Range("E22").Value = 16500
$ws.Range("D23").Value = -31700
$ws.Range("E23").Value = -38400
$ws.Range("D24").Value = -7500
$ws.Range("E24").Value = -7200
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -24200
$ws.Range("E26").Value = -31200
$ws.Range("D27").Value = -24200
$ws.Range("E27").Value = -31200
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -600
$ws.Range("E32").Value = -500
$ws.Range("D33").Value = -24200
$ws.Range("E33").Value = -31200
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -24200
$ws.Range("E35").Value = -31200
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 224900
$ws.Range("E41").Value = 108100
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 54900
$ws.Range("E43").Value = 50100
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 19800
$ws.Range("E45").Value = 14000
$ws.Range("D46").Value = 299600
$ws.Range("E46").Value = 172200
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 2434800
$ws.Range("E48").Value = 2456300
$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 30200
$ws.Range("E52").Value = 26200
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 2764600
$ws.Range("E54").Value = 2654700
$ws.Range("D57").Value = 26800
$ws.Range("E57").Value = 22500
$ws.Range("D58").Value = 96300
$ws.Range("E58").Value = 95100
$ws.Range("D59").Value = 38100
$ws.Range("E59").Value = 43600
$ws.Range("D60").Value = 161200
$ws.Range("E60").Value = 161200
$ws.Range("D61").Value = 1123600
$ws.Range("E61").Value = 989100
$ws.Range("D62").Value = 171800
$ws.Range("E62").Value = 174600
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 1456700
$ws.Range("E66").Value = 1324800
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 549500
$ws.Range("E72").Value = 573700
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 1307900
$ws.Range("E76").Value = 1329900
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = -24200
$ws.Range("E81").Value = -31200
$ws.Range("D83").Value = 27600
$ws.Range("E83").Value = 27600
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = -16600
$ws.Range("E89").Value = 1900
$ws.Range("D91").Value = 0
$ws.Range("E91").Value = 0
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -700
$ws.Range("E94").Value = -2700
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 134100
$ws.Range("E100").Value = 0
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = -100
$ws.Range("D102").Value = 116900
$ws.Range("E102").Value = -1000

# --- 4. Row 58 (Short/Current Long Term Debt): F:J were reported "NA" ----
# instead of carrying forward as 0 once the new quarters were added.
$ws.Range("F58").Value = "NA"
$ws.Range("G58").Value = "NA"
$ws.Range("H58").Value = "NA"
$ws.Range("I58").Value = "NA"
$ws.Range("J58").Value = "NA"

# --- 5. Row 91 (Capital Expenditures): figures were corrected, not just --
# shifted, when this quarter's update was made.
$ws.Range("F91").Value = -100
$ws.Range("G91").Value = 0
$ws.Range("H91").Value = -300
$ws.Range("I91").Value = -800
$ws.Range("J91").Value = -400
$ws.Range("K91").Value = -3600
$ws.Range("L91").Value = -10400
$ws.Range("M91").Value = -25000
